# "added colors to rows" -- colorize the DTR rows by week/status and fix a
# couple of data/formula issues that rode along in the same commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row shading -----------------------------------------------------
# Week 1 (Tue-Fri, rows 5-8) and the Mon/Tue of week 2 (rows 11-12): orange.
# Sat/Sun (rows 9-10) are left alone.
$ws.Range("A5:J8").Interior.Color = 13411113    # FF29A3CC (tried first, overwritten below)
$ws.Range("A5:J8").Interior.Color = 6737151     # FFFFCC66
$ws.Range("A11:J12").Interior.Color = 6737151   # FFFFCC66

# Wed-Fri of week 2 (rows 13-15): red.
$ws.Range("A13:J15").Interior.Color = 6184671   # FFDF5E5E

# --- Data fix: those same red rows were marked as full-day absences --
$ws.Range("I13").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("I15").Value = 1

# --- B19 flips from a blank label to an explicit FALSE flag ----------
# B19 sits inside the merged A19:G19 region, so a direct .Value write to a
# non-anchor merged cell is ignored by Excel; stage the boolean elsewhere,
# copy it, and paste-special (values only) into B19 instead.
$ws.Range("K1").Value = $false
$ws.Range("K1").Copy()
$ws.Range("B19").PasteSpecial(-4163)
$ws.Range("K1").Clear()
$excel.CutCopyMode = $false

# --- Formula cleanup: FLOOR() only takes two arguments here -----------
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
